$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-range cells (shared strings) to their new values.
$ws.Range("C6").Value = "19:45-19:50"
$ws.Range("C7").Value = "19:50-19:55"

# Update the saved selection/active cell on the sheet view from C11 to C10.
$ws.Range("C10").Select()
